# Update result comparison file
# Target sheet: "Comparison_sensitivities" (3rd tab, already the active sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New "InvCost" summary column (R) + ratio helper cells (T/U)
# ---------------------------------------------------------------------
$ws.Range("R1").Value = "InvCost"

$ws.Range("R2").Formula = "=D40"
$ws.Range("R4").Formula = "=`$D`$2"
$ws.Range("R6").Formula = "=D42"

$ws.Range("T4").Formula = "=R2/P4"
$ws.Range("U4").Formula = "=1-T4"
$ws.Range("T5").Formula = "=R6/P4"
$ws.Range("U5").Formula = "=T5-1"

# Match formatting of the neighbouring highlighted summary cells (R4 should
# look like M4:Q4, T4 should look like the highlighted fill without a
# border).
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)

$ws.Range("L4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Borders.LineStyle = -4142

# ---------------------------------------------------------------------
# 2) Four new data rows (sensitivity: investment cost +/-10 %)
# ---------------------------------------------------------------------
$ws.Cells.Item(40,1).Value = "sens_10op_invc_10pdown"
$ws.Range("B40").Formula = "=IF(ISNUMBER(SEARCH(""PV"", A40)),""PV revenue"",""no PV revenue"")"
$ws.Cells.Item(40,3).Value = 239.95405528630451
$ws.Cells.Item(40,4).Value = 1326.412694499294
$ws.Cells.Item(40,5).Value = 261493301.12125921
$ws.Cells.Item(40,6).Value = 15811545.29718912
$ws.Cells.Item(40,7).Value = 176888.9279999978
$ws.Cells.Item(40,8).Value = 32000.007075376489
$ws.Cells.Item(40,9).Value = 9.8181474074492936

$ws.Cells.Item(41,1).Value = "sens_10op_invc_10pdown_PV"
$ws.Range("B41").Formula = "=IF(ISNUMBER(SEARCH(""PV"", A41)),""PV revenue"",""no PV revenue"")"
$ws.Cells.Item(41,3).Value = 186.74954755010739
$ws.Cells.Item(41,4).Value = 1032.309998957538
$ws.Cells.Item(41,5).Value = 261493301.12125921
$ws.Cells.Item(41,6).Value = 6400256.9589656191
$ws.Cells.Item(41,7).Value = 176888.9279999978
$ws.Cells.Item(41,8).Value = 32000.007075376489
$ws.Cells.Item(41,9).Value = 9.8181474074492936

$ws.Cells.Item(42,1).Value = "sens_10op_invc_10pup"
$ws.Range("B42").Formula = "=IF(ISNUMBER(SEARCH(""PV"", A42)),""PV revenue"",""no PV revenue"")"
$ws.Cells.Item(42,3).Value = 272.18328719130898
$ws.Cells.Item(42,4).Value = 1504.5687264186249
$ws.Cells.Item(42,5).Value = 315776480.57557988
$ws.Cells.Item(42,6).Value = 15983677.849763259
$ws.Cells.Item(42,7).Value = 176888.92800000001
$ws.Cells.Item(42,8).Value = 32000.007075376889
$ws.Cells.Item(42,9).Value = 9.8181474074492936

$ws.Cells.Item(43,1).Value = "sens_10op_invc_10pup_PV"
$ws.Range("B43").Formula = "=IF(ISNUMBER(SEARCH(""PV"", A43)),""PV revenue"",""no PV revenue"")"
$ws.Cells.Item(43,3).Value = 218.82991316231221
$ws.Cells.Item(43,4).Value = 1209.6431310916701
$ws.Cells.Item(43,5).Value = 315776480.57557988
$ws.Cells.Item(43,6).Value = 6546056.7125909813
$ws.Cells.Item(43,7).Value = 176888.92800000001
$ws.Cells.Item(43,8).Value = 32000.007075376889
$ws.Cells.Item(43,9).Value = 9.8181474074492936

# ---------------------------------------------------------------------
# 3) Grow the worksheet table to cover the new rows
# ---------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:I43"))

# ---------------------------------------------------------------------
# 4) Restore the selection/scroll state recorded in the saved file
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("W21").Select() | Out-Null
